$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 7
    6  = 2
    7  = 3
    8  = 5
    9  = 5
    10 = 1
    11 = 4
    12 = 4
    13 = 4
    14 = 4
    15 = 1
    16 = 0
    17 = 1
    18 = 3
    19 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
